# Update "想去人数" (want-to-go count) figures that changed between data refreshes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value  = 1305   # was 1304
$wsExhibit.Range("F15").Value = 5643   # was 5642
$wsExhibit.Range("F18").Value = 1126   # was 1125
$wsExhibit.Range("F29").Value = 2941   # was 2939
$wsExhibit.Range("F36").Value = 172    # was 171

# Sheet "本地生活" (local life)
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 1317      # was 1316

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 1317       # was 1316
$wsAll.Range("F6").Value  = 1305       # was 1304
$wsAll.Range("F15").Value = 5643       # was 5642
$wsAll.Range("F18").Value = 1126       # was 1125
$wsAll.Range("F29").Value = 2941       # was 2939
